# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change: cell B11 on sheet "Rules" previously held the text "R40".
# It is now set to hold the text "1" (stored as a new shared string),
# which is also what introduces the new shared-string table entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "1"
